$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the "Value" column (B) for the soil-analysis parameter rows.
$ws.Range("B2").Value = "Soil"
$ws.Range("B3").Value = "test"
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = "Yes"
$ws.Range("B6").Value = 12

$sampleDate = Get-Date -Year 2019 -Month 5 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("B7").Value = $sampleDate
$ws.Range("B7").NumberFormat = "mm-dd-yy"

$ws.Range("B8").Value = 11
$ws.Range("B9").Value = "Silty_Clay"
$ws.Range("B10").Value = "0-30"
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 50
$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 4

# Move the active selection to B8 (matches the saved workbook's view state).
$ws.Range("B8").Select()
